# "Fruta / hortaliza, semanal"
# Insert one new weekly record at row 104 (pushing the existing rows 104-136
# down to 105-137) on the Oregano / Lo Valledor sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a blank row at 104.
$ws.Rows("104:104").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A104").Value = 6
$ws.Range("B104").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C104").Value = "Metropolitana"
$ws.Range("D104").Value = 44559
$ws.Range("E104").Value = 13
$ws.Range("F104").Value = 100112029
$ws.Range("G104").Value = "Orégano"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 33
$ws.Range("K104").Value = 9000
$ws.Range("L104").Value = 10000
$ws.Range("M104").Value = 9455
$ws.Range("N104").Value = "$/docena de atados"
$ws.Range("O104").Value = "Región Metropolitana"
$ws.Range("P104").Value = 3152
$ws.Range("Q104").Value = 3
$ws.Range("R104").Value = "Hortaliza"
